# Add a 2023 column (K) to the Marneuli average monthly remuneration table,
# mirroring the existing 2022 column (J) formatting/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the whole 2022 column (header + data rows) into K to inherit styles,
# then overwrite with the 2023 values.
$ws.Range("J3:J6").Copy($ws.Range("K3:K6"))

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 903.9
$ws.Range("K5").Value = 626.3
$ws.Range("K6").Value = 1102.3
